$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 8-15: the "name" column shifts down the line/extr sequence by two
# positions (two new "line" entries, line7 and line8, are inserted into the
# naming sequence ahead of the extr* entries). Row numbers / A-column ids do
# not change, only the displayed name and the from_bus/to_bus/in_service
# values.
$ws.Range("B8").Value  = "line7"
$ws.Range("C8").Value  = 14
$ws.Range("D8").Value  = 11
$ws.Range("E8").Value  = $true

$ws.Range("B9").Value  = "line8"
$ws.Range("C9").Value  = 16
$ws.Range("D9").Value  = 9
$ws.Range("E9").Value  = $false

$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $true

$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- New rows 16 and 17 (extr7 / extr8), appended after the previous last row.
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false

# Match the bold/bordered/centered style used by the rest of column A.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
